# Daily update at 8 AM UTC
# Appends the next day's row of data to the "Wins Over Time" tracking sheet.
# The last populated row carries a date-only number format (YYYY-MM-DD) to
# visually mark it as the most recent entry; once a new row is appended,
# that marker format must move from the old last row to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the "most recent row" date format (currently on row 67) and the
# regular date/time format used by all the other data rows (e.g. row 2),
# before we change anything.
$lastRowDateFormat = $ws.Range("A67").NumberFormat
$regularDateFormat = $ws.Range("A2").NumberFormat

# Row 67 is no longer the newest row, so it reverts to the regular format.
$ws.Range("A67").NumberFormat = $regularDateFormat

# Append the new day's data as row 68.
$ws.Range("A68").Value = 45655
$ws.Range("B68").Value = 162
$ws.Range("C68").Value = 152
$ws.Range("D68").Value = 157

# Row 68 is now the newest row, so it gets the "most recent" date-only format.
$ws.Range("A68").NumberFormat = $lastRowDateFormat
